# Bug fixes to export_metadata
# This script updates the "gear_type" (column D) sample-type lists on Sheet1,
# fixing a typo (double space / stray space before a semicolon), adding missing
# ice-corer sizes to the "Ice corer ..." list, and filling in a blank cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# 1) Extend the ice-corer list in D19 with the 18 cm and 12 cm corer sizes.
$ws1.Range("D19").Value2 = "Ice corer 18 cm; Ice corer 14 cm; Ice corer 12 cm; Ice corer 9 cm; Suction pump; Slurp gun; Bottle"

# 2) Fix the "Niskin bottle; ... Bottle ; Go-Flo" string (had a stray space
#    before the semicolon) across every cell in column D that used it.
$oldNiskin = "Niskin bottle; Ice corer 14 cm; Ice corer 9 cm; Bucket; Bottle ; Go-Flo"
$newNiskin = "Niskin bottle; Ice corer 14 cm; Ice corer 9 cm; Bucket; Bottle; Go-Flo"

$niskinCells = @("D2","D3","D4","D5","D6","D9","D10","D11","D12","D13","D14","D15","D18","D25","D27","D28","D29","D31","D32","D33")
foreach ($addr in $niskinCells) {
    $cell = $ws1.Range($addr)
    if ($cell.Value2 -eq $oldNiskin) {
        $cell.Value2 = $newNiskin
    }
}

# 3) Fill in the previously empty gear type for row 30 (Phytoplankton taxonomy).
$ws1.Range("D30").Value2 = "Niskin bottle"

# 4) Move the active selection to D27, matching the state the workbook was
#    saved in.
$ws1.Range("D27").Select()
